# Daily attendance processing - 2025-11-18 18:30:39
# Swap the order of the "dnasr281@gmail.com" entry with the name that
# follows it in the "Recorded By" column (column G) so that
# "dnasr281@gmail.com, <other>" becomes "<other>, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -ne $null -and $val -like "dnasr281@gmail.com, *") {
        $rest = $val.Substring("dnasr281@gmail.com, ".Length)
        $cell.Value = "$rest, dnasr281@gmail.com"
    }
}
